$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column values that look like plain numbers stay as text
# (matches original inline-string cell type instead of Excel auto-converting to numeric)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.211.70'
$ws.Range("E2").Value = '  +0.84%  '
$ws.Range("D3").Value = '1.572.51'
$ws.Range("E3").Value = '  +0.80%  '
$ws.Range("E4").Value = '  +0.49%  '
$ws.Range("D5").Value = '211.64'
$ws.Range("E5").Value = '  +2.05%  '
$ws.Range("E6").Value = '  +0.76%  '
$ws.Range("E7").Value = '  +0.42%  '
$ws.Range("D8").Value = '22.08'
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("E9").Value = '  +0.43%  '
$ws.Range("D10").Value = '0.0601'
$ws.Range("E10").Value = '  +0.74%  '
$ws.Range("D11").Value = '0.0870'
$ws.Range("E11").Value = '  +1.08%  '
$ws.Range("D12").Value = '1.793.13'
$ws.Range("E12").Value = '  +0.62%  '
$ws.Range("D13").Value = '1.570.48'
$ws.Range("E13").Value = '  +0.97%  '
$ws.Range("D14").Value = '3.79'
$ws.Range("E14").Value = '  +0.87%  '
$ws.Range("D15").Value = '0.521'
$ws.Range("E15").Value = '  +0.09%  '
$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").Value = '62.38'
$ws.Range("E16").Value = '  +0.47%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '27.169.12'
$ws.Range("E17").Value = '  +0.72%  '
$ws.Range("D18").Value = '7.46'
$ws.Range("E18").Value = '  +1.17%  '
$ws.Range("D19").Value = '216.60'
$ws.Range("E19").Value = '  -0.28%  '
$ws.Range("D20").Value = '0.0₃0703'
$ws.Range("E20").Value = '  -0.41%  '
$ws.Range("E21").Value = '  +0.50%  '
$ws.Range("D22").Value = '4.15'
$ws.Range("E22").Value = '  +1.22%  '
$ws.Range("D23").Value = '9.25'
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("D24").Value = '1.96'
$ws.Range("E24").Value = '  +1.47%  '
$ws.Range("D25").Value = '153.91'
$ws.Range("E25").Value = '  +0.38%  '
$ws.Range("D26").Value = '6.68'
$ws.Range("E26").Value = '  +0.94%  '
$ws.Range("D27").Value = '15.12'
$ws.Range("E27").Value = '  +0.42%  '
$ws.Range("E28").Value = '  +2.46%  '
$ws.Range("E29").Value = '  +0.40%  '
$ws.Range("D31").Value = '0.0474'
$ws.Range("E31").Value = '  +1.20%  '
$ws.Range("D32").Value = '3.25'
$ws.Range("E32").Value = '  +0.36%  '
$ws.Range("D33").Value = '3.19'
$ws.Range("E33").Value = '  +2.37%  '
$ws.Range("D34").Value = '1.454.55'
$ws.Range("E34").Value = '  +2.10%  '
$ws.Range("D35").Value = '1.11'
$ws.Range("E35").Value = '  +5.03%  '
$ws.Range("D36").Value = '1.62'
$ws.Range("E36").Value = '  +0.45%  '
$ws.Range("E37").Value = '  +1.29%  '
$ws.Range("E38").Value = '  +0.95%  '
$ws.Range("E39").Value = '  +1.15%  '
$ws.Range("E40").Value = '  +2.09%  '
$ws.Range("D41").Value = '0.810'
$ws.Range("E41").Value = '  +0.21%  '
$ws.Range("E42").Value = '  +0.50%  '
$ws.Range("D43").Value = '2.35'
$ws.Range("E43").Value = '  +0.76%  '
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  +0.14%  '
$ws.Range("D45").Value = '64.64'
$ws.Range("E45").Value = '  -0.31%  '
$ws.Range("D46").Value = '1.74'
$ws.Range("E46").Value = '  -0.21%  '
$ws.Range("D47").Value = '1.705.13'
$ws.Range("E47").Value = '  +0.53%  '
$ws.Range("D48").Value = '85.91'
$ws.Range("E49").Value = '  +4.06%  '
$ws.Range("D50").Value = '0.0524'
$ws.Range("E50").Value = '  +0.41%  '
$ws.Range("D51").Value = '0.0960'
$ws.Range("E51").Value = '  +0.66%  '
